# Updated symbol list on Sun Feb 12 06:42:06 UTC 2023 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) text values for
# the crypto symbol rows on Sheet1. Values are stored as plain text in the
# source workbook (inline strings), so each write uses a leading apostrophe
# to force Excel to keep the literal text instead of re-interpreting
# numeric- or percent-looking strings as Number cells, then resets the
# cell style back to Normal so no incidental number-format style sticks to
# the cell (matching the original formatting-free cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $ws.Range($addr).Value = "'" + $value
    $ws.Range($addr).Style = "Normal"
}

Set-TextValue "D2" "307.78"
Set-TextValue "E2" "-0.08%"

Set-TextValue "D3" "41.04"
Set-TextValue "E3" "0.17%"

Set-TextValue "E4" "2.08%"

Set-TextValue "D5" "0.07672"
Set-TextValue "E5" "0.75%"

Set-TextValue "D6" "1.645"
Set-TextValue "E6" "1.26%"

Set-TextValue "D7" "0.9160"
Set-TextValue "E7" "1.91%"

Set-TextValue "E8" "0.16%"

Set-TextValue "D9" "0.1245"
Set-TextValue "E9" "14.86%"

Set-TextValue "D10" "0.1830"
Set-TextValue "E10" "3.61%"

Set-TextValue "D11" "0.09195"
Set-TextValue "E11" "-0.38%"

Set-TextValue "D12" "0.04262"
Set-TextValue "E12" "1.41%"

Set-TextValue "E13" "0.04%"

Set-TextValue "D14" "0.001259"
Set-TextValue "E14" "0.51%"

Set-TextValue "D15" "0.005747"
Set-TextValue "E15" "-2.26%"

Set-TextValue "E17" "-0.13%"

Set-TextValue "D18" "4.312"
Set-TextValue "E18" "1.43%"

Set-TextValue "D20" "7.310"
Set-TextValue "E20" "11.22%"

Set-TextValue "E21" "1.53%"

Set-TextValue "E22" "7.99%"

Set-TextValue "D23" "0.04074"
Set-TextValue "E23" "-0.45%"

Set-TextValue "E24" "3.29%"

Set-TextValue "D25" "0.004341"
Set-TextValue "E25" "6.23%"

Set-TextValue "D26" "0.0001272"
Set-TextValue "E26" "-2.22%"

Set-TextValue "D38" "0.02473"
Set-TextValue "E38" "4.27%"

Set-TextValue "D39" "0.05294"
Set-TextValue "E39" "2.29%"

Set-TextValue "D40" "0.007849"
Set-TextValue "E40" "0.87%"

Set-TextValue "E41" "1.21%"

Set-TextValue "D42" "0.006861"
Set-TextValue "E42" "1.35%"

Set-TextValue "D43" "0.001913"
Set-TextValue "E43" "-1.96%"

Set-TextValue "D44" "0.007642"
Set-TextValue "E44" "-10.74%"

Set-TextValue "D45" "0.3063"
Set-TextValue "E45" "-0.33%"

Set-TextValue "D46" "0.00006735"
Set-TextValue "E46" "-2.74%"

Set-TextValue "E47" "0.20%"

Set-TextValue "E48" "1,200.43%"

Set-TextValue "D49" "0.003108"
Set-TextValue "E49" "-26.00%"

Set-TextValue "E50" "0.20%"

Set-TextValue "E51" "0.20%"
